$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-04 Saturday", "2025-10-05 Sunday"),
    @("40×31=1240", "88×62=5456"),
    @("75×49=3675", "81×53=4293"),
    @("22×84=1848", "90×91=8190"),
    @("48×19=912", "49×48=2352"),
    @("82×75=6150", "43×62=2666"),
    @("70×75=5250", "15×96=1440"),
    @("37×98=3626", "62×60=3720"),
    @("34×91=3094", "39×52=2028"),
    @("67×15=1005", "45×69=3105"),
    @("59×62=3658", "44×85=3740"),
    @("94×79=7426", "93×88=8184"),
    @("80×48=3840", "20×76=1520"),
    @("66×23=1518", "54×72=3888"),
    @("17×52=884", "74×62=4588"),
    @("94×14=1316", "92×88=8096"),
    @("89×41=3649", "53×43=2279"),
    @("44×66=2904", "58×82=4756"),
    @("95×65=6175", "19×62=1178"),
    @("37×68=2516", "47×47=2209"),
    @("88×82=7216", "79×75=5925"),
    @("49×47=2303", "16×31=496"),
    @("40×66=2640", "32×15=480"),
    @("30×36=1080", "82×98=8036"),
    @("14×81=1134", "85×20=1700"),
    @("86×32=2752", "27×37=999")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
